$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two ranking swaps).
# Force text number format first so values like "1.020" or "27.516.22" are
# preserved verbatim as text instead of being parsed as numbers/dates.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.516.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.28%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.43%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.55"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.33%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4687"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3980"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.74%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.57"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -10.69%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08036"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.023"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.71"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.58%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.862.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.960"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.63%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.202"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.91%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.14%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.97"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.84%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001036"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.78%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06579"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.42%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.517"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.518.87"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.105.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.60%  "

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.33%  "

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.04"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.65%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.087"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.567"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.28%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.49"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.80%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9592"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.66%  "

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09476"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.71%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.475"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.21%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.602"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.77%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.308"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06101"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.71%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02252"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.92%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.222"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.00%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.136"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.67%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5991"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.79%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1899"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.34"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.61%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.256"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.39%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5697"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.87%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.15"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.72%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.415"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.62%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.946"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06781"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.90%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.07%  "
